$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("products_constraints")

$ws.Range("B1").Value = "p_Names"
$ws.Range("C1").Select()
